$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Age"
$ws.Range("E1").Value = "Role"
$ws.Range("F1").Value = "Address"
$ws.Range("G1").Value = "Company"

$ws.Columns.Item(6).ColumnWidth = 49.166666666666664
